# "semana 17 de 2025"
# Adds a new weekly column (T = week 17) to the IRA/UCI weekly revision
# sheet, mirroring the same "add next epi-week column" edit that created
# columns D..S (weeks 1..16) for previous weeks. Also corrects Q26 (week 14
# for CLINICA LOS ROSALES), which the source data revised from 2 to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell T1 = "17" (week number), styled like the rest of the
# header row (bold + centered, same as D1:S1). Pre-formatting the cell as
# Text before assigning the numeric-looking label keeps it a text value
# (matching t="inlineStr"/shared-string "17") instead of Excel's automatic
# number coercion.
$headerCell = $ws.Cells.Item(1, 20)
$headerCell.NumberFormat = "@"
$headerCell.Font.Bold = $true
$headerCell.HorizontalAlignment = -4108
$headerCell.Value = "17"

# --- Correction: Q26 (week 14, CLINICA LOS ROSALES) changes from 2 to 0.
$ws.Cells.Item(26, 17).Value = 0

# --- New week-17 counts for column T, one per UPGD/row. Rows 18, 31 and 40
# have no reported data for this week (they stay blank, matching the
# source), so they are intentionally skipped.
$week17 = @(
    @{Row=2;  Val=0},
    @{Row=3;  Val=0},
    @{Row=4;  Val=0},
    @{Row=5;  Val=0},
    @{Row=6;  Val=0},
    @{Row=7;  Val=0},
    @{Row=8;  Val=0},
    @{Row=9;  Val=0},
    @{Row=10; Val=0},
    @{Row=11; Val=0},
    @{Row=12; Val=0},
    @{Row=13; Val=0},
    @{Row=14; Val=0},
    @{Row=15; Val=0},
    @{Row=16; Val=0},
    @{Row=17; Val=0},
    @{Row=19; Val=0},
    @{Row=20; Val=0},
    @{Row=21; Val=0},
    @{Row=22; Val=0},
    @{Row=23; Val=0},
    @{Row=24; Val=0},
    @{Row=25; Val=0},
    @{Row=26; Val=10},
    @{Row=27; Val=1},
    @{Row=28; Val=12},
    @{Row=29; Val=0},
    @{Row=30; Val=0},
    @{Row=32; Val=6},
    @{Row=33; Val=0},
    @{Row=34; Val=0},
    @{Row=35; Val=0},
    @{Row=36; Val=0},
    @{Row=37; Val=0},
    @{Row=38; Val=0},
    @{Row=39; Val=0},
    @{Row=41; Val=0},
    @{Row=42; Val=0},
    @{Row=43; Val=0},
    @{Row=44; Val=0},
    @{Row=45; Val=0},
    @{Row=46; Val=0},
    @{Row=47; Val=0},
    @{Row=48; Val=0},
    @{Row=49; Val=1},
    @{Row=50; Val=0},
    @{Row=51; Val=0},
    @{Row=52; Val=0},
    @{Row=53; Val=0}
)

foreach ($entry in $week17) {
    $ws.Cells.Item($entry.Row, 20).Value = $entry.Val
}
